# BOM update: items that were previously marked "N" (not received) have now
# arrived, and one previously-unquantified line item has a received count.
# This mirrors the commit "update BOM, now finished".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 (1u ceramic cap, 1276-1036-1-ND): now received -> Y, keep the old
# "Got 1" note in K and add a new note in L about the rest arriving.
$ws.Range("J6").Value = "Y"
$ws.Range("L6").Value = "now got 11"

# Row 10 (897-1262 camcon connector): quantity received is now known (2),
# so replace the stale placeholder formula with the literal count and mark
# it received.
$ws.Range("I10").Value = 2
$ws.Range("J10").Value = "Y"

# Row 13 (H11856CT-ND camera connector): now received -> Y, add note that
# 2 more have arrived.
$ws.Range("J13").Value = "Y"
$ws.Range("L13").Value = "Have now got 2 more"

# Row 36 (RMCF0603FT56K2CT-ND 56.2k resistor): previously "none arrived",
# now received -> Y, add a note with the updated count.
$ws.Range("J36").Value = "Y"
$ws.Range("L36").Value = "now have 10"

# Leave the selection on the newly-added note cell, matching where the
# author was last working.
$ws.Range("L36").Select() | Out-Null
